$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 / column A: update the exception title text
# (was "Excepção 1 (passo 2) [Pintor já existe no sistema]",
#  now "Excepção 1 (passo 2) [Pintor Existente]")
$ws.Range("A16").Value = "Excepção 1               (passo 2)`n[Pintor Existente]"

# Row 16 / column D keeps its original text
# ("Indica que o pintor já existe no sistema") - no change needed there.

# The active/selected cell when the workbook was last saved moved from A20 to A17
$ws.Range("A17").Select()
